$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column E values for data rows (column deleted from data)
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()

# Row 2 updates
$ws.Range("D2").Value = 0.0343
$ws.Range("F2").Value = -0.00508
$ws.Range("G2").Value = 0.03578012973711164
$ws.Range("H2").Value = 0.03578012973711164
$ws.Range("I2").Value = -0.02383452282301649
$ws.Range("J2").Value = -0.02383452282301649
$ws.Range("K2").Value = -1361
$ws.Range("L2").Value = -0.03097758051667236
$ws.Range("M2").Value = 2486
$ws.Range("N2").Value = 0.09462041456377933
$ws.Range("O2").Value = -1.826598089639971
$ws.Range("P2").Value = 1765
$ws.Range("Q2").Value = 0.06717821066173392
$ws.Range("R2").Value = -1.296840558412932
$ws.Range("S2").Value = 721
$ws.Range("T2").Value = 0.2900241351568785
$ws.Range("U2").Value = 10103
$ws.Range("V2").Value = 0.3845334064110469
$ws.Range("W2").Value = -0.04519042401301591
$ws.Range("X2").Value = 0.09360143399501374
$ws.Range("Y2").Value = -0.1387918580080297
$ws.Range("Z2").Value = 1.29602347956533
$ws.Range("AA2").Value = -0.03089010120286509
$ws.Range("AB2").Value = 0.06836333534660596
$ws.Range("AC2").Value = -0.09925343654947105
$ws.Range("AD2").Value = 12621
$ws.Range("AE2").Value = 535.8488011461467
$ws.Range("AF2").Value = 13156.84880114615
$ws.Range("AG2").Value = 3053.848801146147
$ws.Range("AH2").Value = 0.3336739990533285
$ws.Range("AI2").Value = 0.3064375249597927
$ws.Range("AJ2").Value = 0.1041300812719534
$ws.Range("AK2").Value = 0.09301482897422268
$ws.Range("AL2").Value = 623
$ws.Range("AM2").Value = 623
$ws.Range("AN2").Value = -14.39110604332953
$ws.Range("AO2").Value = -1.651685393258427
$ws.Range("AP2").Value = -3.482153707122174
$ws.Range("AQ2").Value = -1.651685393258427

# Row 3 updates
$ws.Range("D3").Value = 0.0343
$ws.Range("F3").Value = -0.00508
$ws.Range("G3").Value = 0.03578012973711164
$ws.Range("H3").Value = 0.03578012973711164
$ws.Range("I3").Value = -0.02383452282301649
$ws.Range("J3").Value = -0.02383452282301649
$ws.Range("K3").Value = -1361
$ws.Range("L3").Value = -0.03097758051667236
$ws.Range("M3").Value = 2486
$ws.Range("N3").Value = 0.09462041456377933
$ws.Range("O3").Value = -1.826598089639971
$ws.Range("P3").Value = 1765
$ws.Range("Q3").Value = 0.06717821066173392
$ws.Range("R3").Value = -1.296840558412932
$ws.Range("S3").Value = 721
$ws.Range("T3").Value = 0.2900241351568785
$ws.Range("U3").Value = 10103
$ws.Range("V3").Value = 0.3845334064110469
$ws.Range("W3").Value = -0.04519042401301591
$ws.Range("X3").Value = 0.09360143399501374
$ws.Range("Y3").Value = -0.1387918580080297
$ws.Range("Z3").Value = 1.29602347956533
$ws.Range("AA3").Value = -0.03089010120286509
$ws.Range("AB3").Value = 0.06836333534660596
$ws.Range("AC3").Value = -0.09925343654947105
$ws.Range("AD3").Value = 12621
$ws.Range("AE3").Value = 535.8488011461467
$ws.Range("AF3").Value = 13156.84880114615
$ws.Range("AG3").Value = 3053.848801146147
$ws.Range("AH3").Value = 0.3336739990533285
$ws.Range("AI3").Value = 0.3064375249597927
$ws.Range("AJ3").Value = 0.1041300812719534
$ws.Range("AK3").Value = 0.09301482897422268
$ws.Range("AL3").Value = 623
$ws.Range("AM3").Value = 623
$ws.Range("AN3").Value = -14.39110604332953
$ws.Range("AO3").Value = -1.651685393258427
$ws.Range("AP3").Value = -3.482153707122174
$ws.Range("AQ3").Value = -1.651685393258427
